$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) - bold font, centered A1, wrapped B1
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Element Name"
$ws.Range("B1").Value = "Description"

# ---------------------------------------------------------------------------
# Data rows (2-18): element name in column A, description in column B
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Rzepka kolanowa"
$ws.Range("B2").Value = "Jest trzeszczką znajdującą się w mięśnia czworogłowego uda. Jest elementem składowym stawu kolanowgo"

$ws.Range("A3").Value = "Kość piszczelowa"
$ws.Range("B3").Value = "Jest kością o funkcji strukturalnej, tworzącą podudzie. Należy do kości długich. Koniec dalszy kości piszczelowej tworzy kostkę przyśrodkową. "

$ws.Range("A4").Value = "Kość Udowa"
$ws.Range("B4").Value = "To najdłuższa kość ciała ludzkiego. Pokrywa ją gruba warstwa mieśni i w związku z tym jedynie jej niewielki fragment jest dostępny badaniu palpacyjnemu"

$ws.Range("A5").Value = "Więzadło rzepki"
$ws.Range("B5").Value = "Miejsce przyczepu końcowego mięśnia czworogłowego uda pomiędzy rzepką a guzowatością kości piszczelowej"

$ws.Range("A6").Value = "Więzadło krzyżowe przednie"
$ws.Range("B6").Value = "Przebiega od powierzchni przyśrodkowej kłykcia bocznego kości udowej przednio-przyśrodkowo do pola międzykłykciowego przedniego kości piszczelowej"

$ws.Range("A7").Value = "Więzadło krzyżowe tylne"
$ws.Range("B7").Value = "Przebiega od powierzchni przyśrodkowej kłykcia przyśrodkowego do dołu, do pola międzykłykciowego tylnego"

$ws.Range("A8").Value = "Więzadło poprzeczne"
$ws.Range("B8").Value = "Łączy dwie łękotki od przodu, nie jest pokryte błoną maziową (jedyne prawdziwe więzadło śródstawowe ciała ludzkiego)"

$ws.Range("A9").Value = "Łąkotka przyśrodkowa"
$ws.Range("B9").Value = "Półkolista, w kształcie litery C`nPrzyczepia się do pola międzykłykciowego przedniego i tylnego`nMocno przytwierdzona przyśrodkowo do torebki stawowej, a bocznie - do więzadła pobocznego piszczelowego, które ogranicza jej ruchomość"

$ws.Range("A10").Value = "Łąkotka boczna"
$ws.Range("B10").Value = "Prawie całkowicie kolista`nPrzyczpia się do pola międzykłykciowego przedniego i tylnego`nNie jest przytwierdzona do torebki stawowej, więc jest bardziej ruchoma od łąkotki przyśrodkowej"

$ws.Range("A11").Value = "Kłykieć przyśrodkowy"
$ws.Range("A12").Value = "Kłykieć boczny"
$ws.Range("A13").Value = "Kość strzałkowa"
$ws.Range("A14").Value = "Staw piszczelowy"
$ws.Range("A15").Value = "Staw rzepkowy"
$ws.Range("A16").Value = "Staw strzałkowy"

$ws.Range("A17").Value = "Więzadło poboczne strzałkowe"
$ws.Range("B17").Value = "Rozpościera się od bocznej powierzchni kości udowej aż do tzw. głowy kości strzałkowej i odpowiada za stabilność kolana od strony bocznej"

$ws.Range("A18").Value = "Więzadło poboczne piszczelowe"
$ws.Range("B18").Value = "Rozpościera się od przyśrodkowej powierzchni kości udowej aż do przyśrodkowej powierzchni piszczeli, stąd odpowiada za stabilność stawu kolanowego od strony przyśrodkowej (czyli od wewnątrz)."

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# Whole table: column A centered (horizontal+vertical), column B wrapped text
$ws.Range("A1:A18").HorizontalAlignment = -4108
$ws.Range("A1:A18").VerticalAlignment = -4108
$ws.Range("B1:B10").WrapText = $true
$ws.Range("B17:B18").WrapText = $true

# Header row is bold
$ws.Range("A1:B1").Font.Bold = $true

# Row heights for wrapped multi-line rows
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 72
$ws.Rows.Item(10).RowHeight = 57.6
$ws.Rows.Item(17).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 43.2

# Column widths
$ws.Columns.Item(1).ColumnWidth = 26.6
$ws.Columns.Item(2).ColumnWidth = 58.05

# Page setup
$ws.PageSetup.Orientation = 1

# Selection, as left by the author after editing
$ws.Range("B1").Select()
